# Fruta / hortaliza, semanal
# Insert one new weekly record as row 102 in the data table, pushing the
# existing rows 102:163 down to 103:164 (dimension grows from A1:T163 to
# A1:T164).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record by inserting a blank row at position 102.
$ws.Rows("102:102").Insert()

# Populate the newly inserted row with the new price record.
$ws.Range("A102").Value = 7
$ws.Range("B102").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C102").Value = "Ñuble"
$ws.Range("D102").Value = 44518
$ws.Range("E102").Value = 16
$ws.Range("F102").Value = "Fruta"
$ws.Range("G102").Value = 100104
$ws.Range("H102").Value = "Frutos de pepita"
$ws.Range("I102").Value = 100104005
$ws.Range("J102").Value = "Pera"
$ws.Range("K102").Value = "Packham's Triumph"
$ws.Range("L102").Value = "Primera"
$ws.Range("M102").Value = 120
$ws.Range("N102").Value = 10000
$ws.Range("O102").Value = 11000
$ws.Range("P102").Value = 10500
$ws.Range("Q102").Value = "$/caja 16 kilos empedrada"
$ws.Range("R102").Value = "Provincia de Curicó"
$ws.Range("S102").Value = 656
$ws.Range("T102").Value = 16
